# Factor-analysis sheet: drop the 5th ("quartile 4") column from each of the
# Developed / Emerging blocks, so each market block now spans 4 quartile
# columns (0-3) instead of 5 (0-4), and refresh the forward-return row with
# the newly computed values.
#
# Before:  B:F = Developed (quartiles 0-4), G:K = Emerging (quartiles 0-4)
# After :  B:E = Developed (quartiles 0-3), F:I = Emerging (quartiles 0-3)
#
# We do this with real column deletes (shifting everything left) rather than
# Unmerge/Merge: deleting the whole column naturally re-sizes the header
# merge ranges and the sheet dimension/row spans for free, and - unlike the
# Range.Merge() COM call - it does not fabricate new per-edge border styles
# for the cells it touches, so the existing style table is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the old "quartile 4" column of the Developed block (F) - this shifts
# G:K left to F:J and shrinks the B1:F1 merge down to B1:E1 automatically.
$ws.Columns("F:F").Delete()

# Drop the old "quartile 4" column of the Emerging block, which (after the
# shift above) now lives in column J - this shrinks the F1:J1 merge down to
# F1:I1 automatically.
$ws.Columns("J:J").Delete()

# Refresh row 4 (Forward Return) with the newly computed values for the
# remaining 8 columns (B:I).
$ws.Cells.Item(4, 2).Value2 = 0.008672491529713758
$ws.Cells.Item(4, 3).Value2 = 0.006145790245791497
$ws.Cells.Item(4, 4).Value2 = 0.006682893376225255
$ws.Cells.Item(4, 5).Value2 = 0.0076686018312505
$ws.Cells.Item(4, 6).Value2 = 0.01221129846289594
$ws.Cells.Item(4, 7).Value2 = 0.01408019947727526
$ws.Cells.Item(4, 8).Value2 = 0.01184412013195353
$ws.Cells.Item(4, 9).Value2 = 0.0111093832589344
